# Update mods data [2025-12-30 15:11:10]
# Append a new row (50) to the ModCounts sheet with the latest mod count
# for 逃离鸭科夫 on 2025/12/30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new date cell as Text first so that entering "2025/12/30"
# is stored as a literal string instead of being auto-converted into a
# date serial number by Excel's smart data-entry detection.
$ws.Range("A50").NumberFormat = "@"

$ws.Range("A50").Value = "2025/12/30"
$ws.Range("B50").Value = "逃离鸭科夫"
$ws.Range("C50").Value = 1106

# Make the new row match the formatting (centered alignment, default
# General number format) used by the rest of the data rows by copying
# the formatting from the row directly above it. This overwrites the
# temporary Text format applied above while leaving the cell values
# untouched.
$ws.Range("A49:C49").Copy()
$ws.Range("A50:C50").PasteSpecial(-4122)
